$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: shift item-table headers right by two columns (C..G) and
# relabel them; H/I keep GST/Total.
$ws.Range("C1").Value = "ID"
$ws.Range("D1").Value = "DESCRIPTION"
$ws.Range("E1").Value = "QTY"
$ws.Range("F1").Value = "PRICE"
$ws.Range("G1").Value = "TOTAL"

# Clear old line-item block (rows 4-13, cols A-E) before re-writing the
# data in its new location.
$ws.Range("A4:E13").ClearContents()

# Line items, now starting at row 2 and occupying columns C-G.
$items = @(
    @(1, "4Tech keyboard black ", 1, 600, 600),
    @(2, "A4Tech HS-800 headphone ", 1, 900, 900),
    @(3, "Asus Memo Pad Tablet ", 1, 7800, 7800),
    @(4, "HP Desktop C2500 Keyboard+Mouse ", 1, 1500, 1500),
    @(5, "Logitech B170 Wireless Mouse (Black) ", 2, 600, 1200),
    @(6, "Benq G2020HDA Screen. ", 2, 1500, 3000),
    @(7, "Logitech B525 Commercial HD Webcam ", 1, 2000, 2000)
)

$row = 2
foreach ($item in $items) {
    $ws.Cells.Item($row, 3).Value = $item[0]
    $ws.Cells.Item($row, 4).Value = $item[1]
    $ws.Cells.Item($row, 5).Value = $item[2]
    $ws.Cells.Item($row, 6).Value = $item[3]
    $ws.Cells.Item($row, 7).Value = $item[4]
    $row++
}

# Totals block moves from D11:E13 to F9:G11.
$ws.Range("F9").Value = "Sub Total "
$ws.Range("G9").Value = 17000
$ws.Range("F10").Value = "GST 8% "
$ws.Range("G10").Value = 1360
$ws.Range("F11").Value = "Total "
$ws.Range("G11").Value = 18360
